# Admin console CRUD operations: add the "Event" sheet description (mirrors
# the existing "User" sheet) right after "User", and update the active
# selection so that the newly added "Event" sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- Sheet "User": keep its data, just move the selection off of it ---
$userSheet = $wb.Worksheets.Item("User")
$userSheet.Range("E2").Select()

# --- Add the new "Event" sheet right after "User" ---
$eventSheet = $wb.Worksheets.Add($null, $userSheet)
$eventSheet.Name = "Event"

# Make the Event sheet's columns line up with the User sheet visually.
$eventSheet.Columns.Item(1).ColumnWidth = 10.687074829931966

# --- Fill in the Event table: field name / type / size ---
$eventSheet.Range("A2").Value = "data"
$eventSheet.Range("B2").Value = "D"

$eventSheet.Range("A3").Value = "description"
$eventSheet.Range("B3").Value = "A"
$eventSheet.Range("C3").Value = 20

$eventSheet.Range("A4").Value = "userId(mail)"
$eventSheet.Range("B4").Value = "A"
$eventSheet.Range("C4").Value = 100

# Leave the selection on the Event sheet (this also makes it the active tab).
$eventSheet.Range("C4").Select()
